$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 4051
$ws.Range("I74").Value = 3435.3333
$ws.Range("J74").Value = 4666.6665
$ws.Range("K74").Value = 3435.3333
$ws.Range("L74").Value = 4666.6665
$ws.Range("M74").Value = -2499.3333
$ws.Range("N74").Value = -6538.6665
$ws.Range("H76").Value = 3090836.5
$ws.Range("I76").Value = 3371367.2
$ws.Range("J76").Value = 5000
$ws.Range("K76").Value = 3371367.2
$ws.Range("L76").Value = 5000
$ws.Range("M76").Value = -3371052.2
$ws.Range("N76").Value = -5630
$ws.Range("H77").Value = 4051
$ws.Range("I77").Value = 3435.3333
$ws.Range("J77").Value = 4666.6665
$ws.Range("K77").Value = 17176.6665
$ws.Range("L77").Value = 23333.3325
$ws.Range("M77").Value = -12496.6665
$ws.Range("N77").Value = -32693.3325
$ws.Range("H79").Value = 3090836.5
$ws.Range("I79").Value = 3371367.2
$ws.Range("J79").Value = 5000
$ws.Range("K79").Value = 3371367.2
$ws.Range("L79").Value = 5000
$ws.Range("M79").Value = -3370275.2
$ws.Range("N79").Value = -7184
$ws.Range("H86").Value = 2095.1333
$ws.Range("I86").Value = 2047.3
$ws.Range("J86").Value = 2190.8
$ws.Range("K86").Value = 2047.3
$ws.Range("L86").Value = 2190.8
$ws.Range("M86").Value = -924.3
$ws.Range("N86").Value = -4436.8
$ws.Range("H89").Value = 2095.1333
$ws.Range("I89").Value = 2047.3
$ws.Range("J89").Value = 2190.8
$ws.Range("K89").Value = 10236.5
$ws.Range("L89").Value = 10954
$ws.Range("M89").Value = -4620.5
$ws.Range("N89").Value = -22186
$ws.Range("H100").Value = 22225168
$ws.Range("I100").Value = 41668880
$ws.Range("J100").Value = 3783.7144
$ws.Range("K100").Value = 41668880
$ws.Range("L100").Value = 3783.7144
$ws.Range("M100").Value = -41668339
$ws.Range("N100").Value = -4865.7144
$ws.Range("H135").Value = 5145.4
$ws.Range("I135").Value = 3717.111
$ws.Range("J135").Value = 18000
$ws.Range("K135").Value = 33453.999
$ws.Range("L135").Value = 162000
$ws.Range("M135").Value = -30918.999
$ws.Range("N135").Value = -167070
$ws.Range("H137").Value = 814.6
$ws.Range("I137").Value = 809.1539
$ws.Range("J137").Value = 850
$ws.Range("K137").Value = 2427.4617
$ws.Range("L137").Value = 2550
$ws.Range("M137").Value = 122.5383000000002
$ws.Range("N137").Value = -7650
$ws.Range("H138").Value = 3279.2266
$ws.Range("I138").Value = 1360.6222
$ws.Range("J138").Value = 6157.1333
$ws.Range("K138").Value = 4081.8666
$ws.Range("L138").Value = 18471.3999
$ws.Range("M138").Value = 1058.1334
$ws.Range("N138").Value = -28751.3999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 335016.75
$ws.Range("I32").Value = 2259.4075
$ws.Range("J32").Value = 5725685.5
$ws.Range("K32").Value = 2259.4075
$ws.Range("L32").Value = 5725685.5
$ws.Range("M32").Value = -1972.4075
$ws.Range("N32").Value = -5726259.5
$ws.Range("H123").Value = 28000
$ws.Range("J123").Value = 28000
$ws.Range("L123").Value = 28000
$ws.Range("N123").Value = -37800

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2879.2354
$ws.Range("I20").Value = 2090.3635
$ws.Range("J20").Value = 4325.5
$ws.Range("K20").Value = 2090.3635
$ws.Range("L20").Value = 4325.5
$ws.Range("M20").Value = -1843.3635
$ws.Range("N20").Value = -4819.5
$ws.Range("H134").Value = 12299.692
$ws.Range("I134").Value = 733.55554
$ws.Range("J134").Value = 38323.5
$ws.Range("K134").Value = 2200.66662
$ws.Range("L134").Value = 114970.5
$ws.Range("M134").Value = 334.33338
$ws.Range("N134").Value = -120040.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2163.4043
$ws.Range("I31").Value = 1740.7222
$ws.Range("J31").Value = 3546.7273
$ws.Range("K31").Value = 1740.7222
$ws.Range("L31").Value = 3546.7273
$ws.Range("M31").Value = -1445.7222
$ws.Range("N31").Value = -4136.7273
$ws.Range("H34").Value = 2163.4043
$ws.Range("I34").Value = 1740.7222
$ws.Range("J34").Value = 3546.7273
$ws.Range("K34").Value = 1740.7222
$ws.Range("L34").Value = 3546.7273
$ws.Range("M34").Value = -1538.7222
$ws.Range("N34").Value = -3950.7273
$ws.Range("H58").Value = 2748588
$ws.Range("I58").Value = 1062.2307
$ws.Range("J58").Value = 5496114
$ws.Range("K58").Value = 1062.2307
$ws.Range("L58").Value = 5496114
$ws.Range("M58").Value = -859.2307000000001
$ws.Range("N58").Value = -5496520
$ws.Range("H132").Value = 24346.373
$ws.Range("I132").Value = 621.0294
$ws.Range("J132").Value = 113975.445
$ws.Range("K132").Value = 1863.0882
$ws.Range("L132").Value = 341926.335
$ws.Range("M132").Value = 666.9117999999999
$ws.Range("N132").Value = -346986.335
$ws.Range("H134").Value = 37079.25
$ws.Range("I134").Value = 39880.348
$ws.Range("J134").Value = 665
$ws.Range("K134").Value = 119641.044
$ws.Range("L134").Value = 1995
$ws.Range("M134").Value = -117106.044
$ws.Range("N134").Value = -7065
$ws.Range("H136").Value = 2748588
$ws.Range("I136").Value = 1062.2307
$ws.Range("J136").Value = 5496114
$ws.Range("K136").Value = 3186.6921
$ws.Range("L136").Value = 16488342
$ws.Range("M136").Value = -636.6921000000002
$ws.Range("N136").Value = -16493442

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 11364816
$ws.Range("I131").Value = 1425.7142
$ws.Range("J131").Value = 13514647
$ws.Range("K131").Value = 4277.142599999999
$ws.Range("L131").Value = 40543941
$ws.Range("M131").Value = 762.8574000000008
$ws.Range("N131").Value = -40554021

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 85575.25
$ws.Range("I132").Value = 68519.664
$ws.Range("J132").Value = 114001.22
$ws.Range("K132").Value = 205558.992
$ws.Range("L132").Value = 342003.66
$ws.Range("M132").Value = -203028.992
$ws.Range("N132").Value = -347063.66

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 726.5111000000001
$ws.Range("I22").Value = 510.7143
$ws.Range("J22").Value = 1081.9412
$ws.Range("K22").Value = 510.7143
$ws.Range("L22").Value = 1081.9412
$ws.Range("M22").Value = -215.7143
$ws.Range("N22").Value = -1671.9412
$ws.Range("H27").Value = 726.5111000000001
$ws.Range("I27").Value = 510.7143
$ws.Range("J27").Value = 1081.9412
$ws.Range("K27").Value = 510.7143
$ws.Range("L27").Value = 1081.9412
$ws.Range("M27").Value = -403.7143
$ws.Range("N27").Value = -1295.9412
$ws.Range("H132").Value = 1354626
$ws.Range("I132").Value = 3712.25
$ws.Range("J132").Value = 5557469
$ws.Range("K132").Value = 11136.75
$ws.Range("L132").Value = 16672407
$ws.Range("M132").Value = -8606.75
$ws.Range("N132").Value = -16677467

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 19684624
$ws.Range("I132").Value = 31251262
$ws.Range("J132").Value = 1178002.2
$ws.Range("K132").Value = 93753786
$ws.Range("L132").Value = 3534006.6
$ws.Range("M132").Value = -93751256
$ws.Range("N132").Value = -3539066.6
$ws.Range("H136").Value = 4314845
$ws.Range("I136").Value = 6052.737
$ws.Range("J136").Value = 12501550
$ws.Range("K136").Value = 18158.211
$ws.Range("L136").Value = 37504650
$ws.Range("M136").Value = -15608.211
$ws.Range("N136").Value = -37509750
